# Fruta / hortaliza, semanal
#
# The edit re-shuffles the per-row "weekly observation" data (date, quality,
# volume, min/max/weighted prices, commercialisation unit, origin, $/kg and
# kg/unit) across the existing data rows (2-23, excluding 4 and 5 which are
# unchanged). Column A/B/C/E/F/G/H/I/J/K (market/product identity columns)
# stay put; only D, L, M, N, O, P, Q, R, S, T move between rows.
#
# targetRow -> sourceRow : targetRow ends up with the pre-edit values that
# sourceRow held before any writes happened.
$moveMap = @{
    2  = 16
    3  = 9
    6  = 3
    7  = 10
    8  = 11
    9  = 2
    10 = 14
    11 = 15
    12 = 20
    13 = 23
    14 = 21
    15 = 12
    16 = 13
    17 = 6
    18 = 7
    19 = 17
    20 = 18
    21 = 19
    22 = 8
    23 = 22
}

$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot every row involved BEFORE any writes, since several rows both
# donate data to one row and receive data from another (a permutation/cycle),
# so writes must not clobber values still needed later.
$snapshot = @{}
$rowsInvolved = @(2,3,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23)
foreach ($r in $rowsInvolved) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Apply the permutation.
foreach ($targetRow in $rowsInvolved) {
    $sourceRow = $moveMap[$targetRow]
    $sourceData = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $sourceData[$c]
    }
}
